$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "time_taken"
$ws.Range("F1").Style = $ws.Range("E1").Style

$ws.Range("F2").Value = "2021-10-05 10:51:51.097951"
$ws.Range("F3").Value = "2021-10-05 10:51:51.097970"
$ws.Range("F4").Value = "2021-10-05 10:51:51.097977"
$ws.Range("F5").Value = "2021-10-05 10:51:51.097983"
